# Update "想去人数" (number of people wanting to go) figures for a few
# events on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 7830
$wsExhibit.Range("F10").Value = 454
$wsExhibit.Range("F17").Value = 5705
$wsExhibit.Range("F18").Value = 164
$wsExhibit.Range("F19").Value = 235
$wsExhibit.Range("F20").Value = 1345
$wsExhibit.Range("F21").Value = 230
$wsExhibit.Range("F22").Value = 350

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7830
$wsAll.Range("F10").Value = 454
$wsAll.Range("F18").Value = 5705
$wsAll.Range("F20").Value = 164
$wsAll.Range("F21").Value = 235
$wsAll.Range("F22").Value = 1345
$wsAll.Range("F23").Value = 230
$wsAll.Range("F24").Value = 350
